$d = $word.ActiveDocument
Write-Host "Initial paragraphs:" $d.Paragraphs.Count

# Delete paragraphs after the 2nd one (paragraphs 3..6), from the end backwards
for ($i = $d.Paragraphs.Count; $i -ge 3; $i--) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Delete()
}

# Delete paragraph 1 entirely (its text + paragraph mark)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Delete()
Write-Host "After removing paragraph1:" $d.Paragraphs.Count
Write-Host "Content: [" $d.Content.Text "]"

# Now remove the remaining text in paragraph 1 ("Used 3- decoders...")
$p = $d.Paragraphs.Item(1)
$rng = $p.Range
Write-Host "p.Range Start=" $rng.Start "End=" $rng.End
$textRng = $d.Range($rng.Start, $rng.End - 1)
Write-Host "textRng: [" $textRng.Text "]"
$textRng.Delete()
Write-Host "After removing text:" $d.Paragraphs.Count
Write-Host "Content: [" $d.Content.Text "]"
